# Update column G ("K") values on Sheet1 rows 2-23 to reflect
# regenerated save_data (K replaces Strike#, std/mean recalculated,
# s_vals recalculated and written).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 4
    4  = 2
    5  = 2
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 2
    21 = 3
    22 = 1
    23 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
